$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = ''''
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = '''1.0'
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = '''0.00'
$ws.Range("C9").Value = 77
$ws.Range("G9").Value = '''36344.00'
$ws.Range("A10").Value = 'P. point'
$ws.Range("C10").Value = 87
$ws.Range("D10").Value = '''6'
$ws.Range("E10").Value = 'On board'
$ws.Range("F10").Value = 136
$ws.Range("G10").Value = '''11832.00'
$ws.Range("A11").Value = 'Each'
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = '''3.0'
$ws.Range("E11").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 23
$ws.Range("G11").Value = '''184.00'
$ws.Range("D12").Value = '''4.0'
$ws.Range("E12").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 50
$ws.Range("G12").Value = '''2850.00'
$ws.Range("C13").Value = 23
$ws.Range("D13").Value = '''7.0'
$ws.Range("E13").Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = '''690.00'
$ws.Range("C14").Value = 95
$ws.Range("D14").Value = '''10.0'
$ws.Range("E14").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F14").Value = 303
$ws.Range("G14").Value = '''28785.00'
$ws.Range("A15").Value = ''''
$ws.Range("C15").Value = 60
$ws.Range("D15").Value = '''11.0'
$ws.Range("E15").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = '''0.00'
$ws.Range("A16").Value = 'R. mtr.'
$ws.Range("C16").Value = 73
$ws.Range("D16").Value = '''17'
$ws.Range("E16").Value = '25 mm'
$ws.Range("F16").Value = 56
$ws.Range("G16").Value = '''4088.00'
$ws.Range("A17").Value = ''''
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = '''12.0'
$ws.Range("E17").Value = 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = '''0.00'
$ws.Range("A18").Value = 'Mtr.'
$ws.Range("C18").Value = 45
$ws.Range("D18").Value = '''19'
$ws.Range("E18").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F18").Value = 81
$ws.Range("G18").Value = '''3645.00'
$ws.Range("A19").Value = 'Mtr.'
$ws.Range("C19").Value = 79
$ws.Range("D19").Value = '''20'
$ws.Range("E19").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F19").Value = 122
$ws.Range("G19").Value = '''9638.00'
$ws.Range("A20").Value = 'Set'
$ws.Range("C20").Value = 52
$ws.Range("D20").Value = '''13.0'
$ws.Range("E20").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F20").Value = 5733
$ws.Range("G20").Value = '''298116.00'
$ws.Range("A21").Value = ''''
$ws.Range("C21").Value = 66
$ws.Range("D21").Value = '''14.0'
$ws.Range("E21").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = '''0.00'
$ws.Range("A22").Value = 'Mtr.'
$ws.Range("C22").Value = 86
$ws.Range("D22").Value = '''23'
$ws.Range("E22").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = '''1720.00'
$ws.Range("A23").Value = 'Each'
$ws.Range("C23").Value = 70
$ws.Range("D23").Value = '''25'
$ws.Range("E23").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F23").Value = 1890
$ws.Range("G23").Value = '''132300.00'
$ws.Range("A24").Value = ''''
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = '''16.0'
$ws.Range("E24").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = '''0.00'
$ws.Range("A25").Value = 'Each'
$ws.Range("C25").Value = 95
$ws.Range("D25").Value = '''27'
$ws.Range("E25").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F25").Value = 492
$ws.Range("G25").Value = '''46740.00'
$ws.Range("A26").Value = ''''
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = '''17.0'
$ws.Range("E26").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = '''0.00'
$ws.Range("A27").Value = ''''
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = '''31'
$ws.Range("E27").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = '''0.00'
$ws.Range("C28").Value = 53
$ws.Range("D28").Value = '''18.0'
$ws.Range("E28").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("C29").Value = 54
$ws.Range("A30").Value = ''''
$ws.Range("C30").Value = 18
$ws.Range("D30").Value = '''36'
$ws.Range("E30").Value = 'Total'
$ws.Range("C31").Value = 18
$ws.Range("G33").Value = '''576932.00'
$ws.Range("H33").Value = 576932.00
$ws.Range("G35").Value = '''576932.00'
$ws.Range("H35").Value = 576932.00
